$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date values in column C (rows 2-6), shifting the original
# 44154.875 date forward so each row gets a distinct day.
$ws.Range("C2").Value = 44184.875
$ws.Range("C3").Value = 44185.875
$ws.Range("C4").Value = 44186.875
$ws.Range("C5").Value = 44187.875
$ws.Range("C6").Value = 44188.875

# Update the view: clear the frozen/scrolled topLeftCell and move the
# active selection to D10.
$ws.Range("D10").Select()
